$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.421.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +2.92%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.064.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +1.86%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.14%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''549.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.83%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''140.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +2.36%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.04%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.061.97'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +2.02%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.98%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''6.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +7.08%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +1.16%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.455'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.93%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000228'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.96%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''34.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +1.27%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.567.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.29%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''63.575.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +3.14%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.073.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +2.35%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -1.31%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''6.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.49%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''485.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +3.18%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +3.16%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.677'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.55%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +3.93%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''80.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.01%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''12.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +4.33%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.22%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +2.56%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''7.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.28%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''2.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +5.61%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -0.01%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''26.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +1.83%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +0.46%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''2.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +6.64%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''5.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +3.34%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''55.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.14%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''5.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.90%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''467.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.97%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.0822'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +3.42%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0397'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.33%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''3.063.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -3.80%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +1.12%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''8.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.91%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''2.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.50%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''28.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +2.23%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.256'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +3.46%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -0.09%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.91%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  +1.35%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.0₃0512'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +2.20%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''116.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.80%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''2.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +2.46%  '
$ws.Range("E51").Style = "Normal"
